$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (standard VBA RGB encoding: r + g*256 + b*65536)
$red    = 6184671   # FFDF5E5E - absent / no time in-out
$orange = 6737151   # FFFFCC66 - half day / undertime

# Row 8  (02-13-2015, Friday)  -> absent, highlight red, mark 1 "no of hours late"-ish col I
$ws.Range("A8:J8").Interior.Color = $red
$ws.Range("I8").Value = 1

# Row 11 (02-16-2015, Monday) -> absent, highlight red
$ws.Range("A11:J11").Interior.Color = $red
$ws.Range("I11").Value = 1

# Row 14 (02-19-2015, Thursday) -> absent, highlight red
$ws.Range("A14:J14").Interior.Color = $red
$ws.Range("I14").Value = 1

# Row 15 (02-20-2015, Friday) -> half day, highlight orange
$ws.Range("A15:J15").Interior.Color = $orange

# B19 used to hold a blank-space text placeholder; now a boolean FALSE
$ws.Range("B19").Value = $false

# Fix FLOOR() calls that were called with an extra (unsupported) third argument
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'

Write-Host "done"
